$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.520.19"
$ws.Range("E2").Value = "  +0.05%  "

$ws.Range("D3").Value = "1.822.85"
$ws.Range("E3").Value = "  -0.24%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.38"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.43%  "

$ws.Range("E6").Value = "  +0.08%  "

$ws.Range("E7").Value = "  -5.46%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3947"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.39%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08213"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +7.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.111"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.86%  "

$ws.Range("E11").Value = "  -0.55%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.326"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.12%  "

$ws.Range("E13").Value = "  +0.35%  "

$ws.Range("E14").Value = "  +0.13%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.515"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.84%  "

$ws.Range("D16").Value = "1.824.10"
$ws.Range("E16").Value = "  +0.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.00001128"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.40%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "92.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.03%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06658"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.84%  "

$ws.Range("E20").Value = "  +0.45%  "

$ws.Range("E21").Value = "  +0.05%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.091"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.21%  "

$ws.Range("D23").Value = "28.548.36"
$ws.Range("E23").Value = "  +0.14%  "

$ws.Range("E24").Value = "  +2.25%  "

$ws.Range("E25").Value = "  -0.39%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "21.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +3.29%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.52"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.87%  "

$ws.Range("D28").Value = "2.034.64"
$ws.Range("E28").Value = "  +0.04%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.402"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.19%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.08"
$ws.Range("D30").Style = "Normal"

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.114"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1092"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.760"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.33%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.658"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.53%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07064"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.56%  "

$ws.Range("E36").Value = "  -0.07%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.281"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.33%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02351"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.23%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.840"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.23%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.6319"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.32%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "11.29"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.74%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.183"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.44%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.001"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.08%  "

$ws.Range("E44").Value = "  +0.04%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.51"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.03%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5932"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.08%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.729"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.09%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.990"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.71%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.188"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06904"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.28%  "
